# 自动更新Excel文件
# For every data row (row 2..99) in column E ("剩余" / remaining days):
#   - if remaining == 1, the cycle restarts: remaining becomes 10 and the
#     start date in column F ("开始时间", stored as an integer YYYYMMDD) is
#     advanced by 10 days.
#   - otherwise remaining is simply decremented by 1 and the date is left
#     untouched.
# Rows whose F value is not a well-formed 8-digit YYYYMMDD number (e.g. a
# corrupted date) are left completely unmodified, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$daysInMonthTable = @(31,28,31,30,31,30,31,31,30,31,30,31)

$firstDataRow = 2
$lastDataRow = 99

for ($rowNum = $firstDataRow; $rowNum -le $lastDataRow; $rowNum++) {

    $remainingCell = $ws.Cells.Item($rowNum, 5)
    $startDateCell = $ws.Cells.Item($rowNum, 6)

    $remainingValue = [int]$remainingCell.Value()
    $startDateValue = [int]$startDateCell.Value()

    $startDateStr = [string]$startDateValue
    $dateIsWellFormed = $false
    if ($startDateStr.Length -eq 8) {
        $dateIsWellFormed = $true
    }

    if ($dateIsWellFormed) {
        if ($remainingValue -eq 1) {
            $newRemainingValue = 10

            $curYear = [int]([math]::Floor($startDateValue / 10000))
            $curMonth = [int]([math]::Floor(($startDateValue % 10000) / 100))
            $curDay = [int]($startDateValue % 100)

            $curDay = $curDay + 10

            $safetyCounter = 0
            while ($safetyCounter -lt 24) {
                $isLeapYear = $false
                if ((($curYear % 4) -eq 0) -and ((($curYear % 100) -ne 0) -or (($curYear % 400) -eq 0))) {
                    $isLeapYear = $true
                }

                $daysInCurMonth = $daysInMonthTable[$curMonth - 1]
                if (($curMonth -eq 2) -and $isLeapYear) {
                    $daysInCurMonth = 29
                }

                if ($curDay -gt $daysInCurMonth) {
                    $curDay = $curDay - $daysInCurMonth
                    $curMonth = $curMonth + 1
                    if ($curMonth -gt 12) {
                        $curMonth = 1
                        $curYear = $curYear + 1
                    }
                } else {
                    break
                }
                $safetyCounter = $safetyCounter + 1
            }

            $newStartDateValue = ($curYear * 10000) + ($curMonth * 100) + $curDay

            $remainingCell.Value = $newRemainingValue
            $startDateCell.Value = $newStartDateValue
        } else {
            $newRemainingValue = $remainingValue - 1
            $remainingCell.Value = $newRemainingValue
        }
    }
}
